$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 76
$ws.Range("G2").Value = 45
$ws.Range("G3").Value = 55
$ws.Range("F4").Value = 4630
$ws.Range("F5").Value = 1839
$ws.Range("F6").Value = 133
$ws.Range("F8").Value = 3110
$ws.Range("F11").Value = 262
$ws.Range("F12").Value = 626
$ws.Range("F17").Value = 1775
$ws.Range("F18").Value = 1328
$ws.Range("F20").Value = 1600
$ws.Range("F21").Value = 127
$ws.Range("F23").Value = 5
$ws.Range("F28").Value = 98
$ws.Range("F32").Value = 3810
$ws.Range("F33").Value = 759
$ws.Range("F34").Value = 68
$ws.Range("F35").Value = 738
$ws.Range("F36").Value = 58
$ws.Range("F37").Value = 1810

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 23

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 76
$ws.Range("G2").Value = 45
$ws.Range("G3").Value = 55
$ws.Range("F4").Value = 4630
$ws.Range("F5").Value = 1839
$ws.Range("F6").Value = 133
$ws.Range("F8").Value = 3110
$ws.Range("F11").Value = 262
$ws.Range("F12").Value = 626
$ws.Range("F15").Value = 23
$ws.Range("F18").Value = 1775
$ws.Range("F21").Value = 1600
$ws.Range("F22").Value = 127
$ws.Range("F24").Value = 5
$ws.Range("F29").Value = 98
$ws.Range("F33").Value = 3810
$ws.Range("F35").Value = 759
$ws.Range("F36").Value = 68
$ws.Range("F37").Value = 738
$ws.Range("F38").Value = 58
$ws.Range("F39").Value = 1810
